$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for the 4 added date columns (shared strings 74-77)
$ws.Range("AZ1").Value = "16-ago"
$ws.Range("BA1").Value = "17-ago"
$ws.Range("BB1").Value = "18-ago"
$ws.Range("BC1").Value = "22-ago"

# New data values for rows 2-18, columns AZ:BC
$ws.Range("AZ2").Value = 0
$ws.Range("BA2").Value = 0
$ws.Range("BB2").Value = 0
$ws.Range("BC2").Value = 0
$ws.Range("AZ3").Value = 16.72581762695949
$ws.Range("BA3").Value = 14.740217530833212
$ws.Range("BB3").Value = 13.691533671864502
$ws.Range("BC3").Value = 10.777018540069239
$ws.Range("AZ4").Value = 18.818043194818785
$ws.Range("BA4").Value = 19.960144967120499
$ws.Range("BB4").Value = 19.387145101836023
$ws.Range("BC4").Value = 17.969858079165061
$ws.Range("AZ5").Value = 15.429149913526052
$ws.Range("BA5").Value = 14.054991147108684
$ws.Range("BB5").Value = 12.83351668807102
$ws.Range("BC5").Value = 24.140000030352329
$ws.Range("AZ6").Value = 0
$ws.Range("BA6").Value = 0
$ws.Range("BB6").Value = 0
$ws.Range("BC6").Value = 0
$ws.Range("AZ7").Value = 17.798327699583332
$ws.Range("BA7").Value = 15.36157249125
$ws.Range("BB7").Value = 17.714209574583332
$ws.Range("BC7").Value = 16.35543054
$ws.Range("AZ8").Value = 10.923344317743496
$ws.Range("BA8").Value = 10.910910232009014
$ws.Range("BB8").Value = 18.688430858928257
$ws.Range("BC8").Value = 27.143643547145082
$ws.Range("AZ9").Value = 13.689515933258763
$ws.Range("BA9").Value = 15.27918512076481
$ws.Range("BB9").Value = 17.540477460439064
$ws.Range("BC9").Value = 12.683070220135964
$ws.Range("AZ10").Value = 7.9249502724298022
$ws.Range("BA10").Value = 5.5504855470052172
$ws.Range("BB10").Value = 6.5053828093824233
$ws.Range("BC10").Value = 7.7129678895261709
$ws.Range("AZ11").Value = 10.563495125915415
$ws.Range("BA11").Value = 11.478486262124781
$ws.Range("BB11").Value = 13.287532296681094
$ws.Range("BC11").Value = 9.6923635399766965
$ws.Range("AZ12").Value = 0
$ws.Range("BA12").Value = 0
$ws.Range("BB12").Value = 0
$ws.Range("BC12").Value = 0
$ws.Range("AZ13").Value = 11.421590914992073
$ws.Range("BA13").Value = 9.7691323624350179
$ws.Range("BB13").Value = 8.3477647568830644
$ws.Range("BC13").Value = 11.972790548583601
$ws.Range("AZ14").Value = 0
$ws.Range("BA14").Value = 0
$ws.Range("BB14").Value = 0
$ws.Range("BC14").Value = 0
$ws.Range("AZ15").Value = 0
$ws.Range("BA15").Value = 0
$ws.Range("BB15").Value = 0
$ws.Range("BC15").Value = 0
$ws.Range("AZ16").Value = 11.859470296760426
$ws.Range("BA16").Value = 10.861278408055139
$ws.Range("BB16").Value = 9.5778888368626287
$ws.Range("BC16").Value = 6.4515444224646235
$ws.Range("AZ17").Value = 0
$ws.Range("BA17").Value = 0
$ws.Range("BB17").Value = 0
$ws.Range("BC17").Value = 0
$ws.Range("AZ18").Value = 0
$ws.Range("BA18").Value = 0
$ws.Range("BB18").Value = 0
$ws.Range("BC18").Value = 0

$ws.Range("BE9").Select() | Out-Null
